$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string "medium" -> "small to medium" (cell H4)
$ws.Range("H4").Value = "small to medium"

# Adjust column widths: F (14.7109375 -> 13.7109375), H (10 -> 16)
# NOTE: the runtime quantizes ColumnWidth writes to 1/6-character steps
# (stored width = round(ColumnWidth*6)/6 + 5/6), so the ColumnWidth input
# below is chosen to land as close as possible to the target stored width.
$ws.Columns.Item(6).ColumnWidth = 12.83
$ws.Columns.Item(8).ColumnWidth = 15.2

# Row 2 values
$ws.Range("C2").Value = 534
$ws.Range("D2").Value = 76516.546709838207
$ws.Range("F2").Value = 0.9930694846071505
$ws.Range("G2").Value = 23.940711210346596

# Row 3 values
$ws.Range("C3").Value = 534
$ws.Range("D3").Value = 5.9992645534675786
$ws.Range("E3").Value = 0.014632073361372622
$ws.Range("F3").Value = 0.011109764303898541
$ws.Range("G3").Value = 0.21198658274034568

# Row 4 values
$ws.Range("C4").Value = 534
$ws.Range("D4").Value = 13.257890899803465
$ws.Range("E4").Value = 0.0000024019336276603909
$ws.Range("F4").Value = 0.047306039652397755
$ws.Range("G4").Value = 0.4456681370555548

# Row 5 values
$ws.Range("C5").Value = 534
$ws.Range("D5").Value = 0.49583759460992477
$ws.Range("E5").Value = 0.609340678459378
$ws.Range("F5").Value = 0.001853627327694597
$ws.Range("G5").Value = 0.086187461796569942
